$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: ResellerTests in column A, Y in column C (B left blank)
$ws.Range("A6").Value = "ResellerTests"
$ws.Range("C6").Value = "Y"

# Move the active selection to reflect the new last empty row
$ws.Range("A7").Select()
